$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('old times', '昔|むかし'),
    @('old tale', '昔話|むかしばなし'),
    @('ancient times', '大昔|おおむかし'),
    @('once upon a time', '昔々|むかしむかし'),
    @('people', '人々|ひとびと'),
    @('sometimes', '時々|ときどき'),
    @('various', '色々な|いろいろな'),
    @('God', '神さま|かみさま'),
    @('shrine', '神社|じんじゃ'),
    @('Shinto', '神道|しんとう'),
    @('Kobe City', '神戸市|こうべし'),
    @('early', '早い|はやい'),
    @('to get up early', '早起きする|はやおきする'),
    @('early morning', '早朝|そうちょう'),
    @('to get up', '起きる|おきる'),
    @('to wake someone up', '起こす|おこす'),
    @('reboot', '再起動|さいきどう'),
    @('cow', '牛|うし'),
    @('milk', '牛乳|ぎゅうにゅう'),
    @('beef', '牛肉|ぎゅうにく'),
    @('calf; veal', '子牛|こうし'),
    @('to use', '使う|つかう'),
    @('ambassador', '大使|たいし'),
    @('"Occupied"', '使用中|しようちゅう'),
    @('errand', 'お使い|おつかい'),
    @('to work', '働く|はたらく'),
    @('both husband and wife working', '共働き|ともばたらき'),
    @('labor', '労働|ろうどう'),
    @('to bring (a person) back', '連れて帰る|つれてかえる'),
    @('United Nations', '国連|こくれん'),
    @('consecutive holidays', '連休|れんきゅう'),
    @('to separate', '別れる|わかれる'),
    @('not in particular', '別に|べつに'),
    @('special', '特別な|とくべつな'),
    @('discrimination', '差別|さべつ'),
    @('separately', '別々に|べつべつに'),
    @('once', '一度|いちど'),
    @('near future', '今度|こんど'),
    @('temperature', '温度|おんど'),
    @('30 degrees', '三十度|さんじゅうど'),
    @('attitude', '態度|たいど'),
    @('red color', '赤|あか'),
    @('red', '赤い|あかい'),
    @('baby', '赤ちゃん|あかちゃん'),
    @('the equator', '赤道|せきどう'),
    @('the Red Cross', '赤十字|せきじゅうじ'),
    @('blue color', '青|あお'),
    @('blue', '青い|あおい'),
    @('youth', '青年|せいねん'),
    @('blue sky', '青空|あおぞら'),
    @('green light', '青信号|あおしんごう'),
    @('color', '色|いろ'),
    @('various', '色々な|いろいろな'),
    @('scenery', '景色|けしき'),
    @('characteristic', '特色|とくしょく'),
)

$startRow = 85
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
